$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are kept as exact text (preserve trailing zeros / formatting)
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell content updates per the data refresh
$ws.Range("D2").Value = "29.897.79"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.889.57"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "0.7684"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").Value = "242.55"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.3137"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").Value = "25.66"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "0.07132"
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("D11").Value = "0.08530"
$ws.Range("E11").Value = "  +4.94%  "
$ws.Range("D12").Value = "0.7642"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "1.906.70"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "5.371"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "93.73"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "6.143"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").Value = "29.911.36"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "13.77"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").Value = "244.37"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "0.000007819"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").Value = "0.9992"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "8.020"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "0.1629"
$ws.Range("E24").Value = "  +2.65%  "
$ws.Range("D25").Value = "9.385"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "163.19"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "2.040"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "1.518"
$ws.Range("E29").Value = "  +4.15%  "
$ws.Range("D30").Value = "1.534"
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").Value = "4.509"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "4.119"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").Value = "0.05453"
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("D34").Value = "1.244"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").Value = "0.7459"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "2.703"
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("D38").Value = "0.01949"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "0.4475"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").Value = "1.102.20"
$ws.Range("E41").Value = "  -3.75%  "
$ws.Range("D42").Value = "73.22"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").Value = "6.084"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").Value = "0.8550"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "103.10"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.688"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.871"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").Value = "3.038"
$ws.Range("E49").Value = "  -2.75%  "
$ws.Range("D50").Value = "2.034.18"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "0.06086"
$ws.Range("E51").Value = "  +0.39%  "
